$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 4915.3335
$ws.Range("I51").Value = 4164.3335
$ws.Range("J51").Value = 5666.3335
$ws.Range("K51").Value = 4164.3335
$ws.Range("L51").Value = 5666.3335
$ws.Range("M51").Value = -3680.3335
$ws.Range("N51").Value = -6634.3335
# Row 98
$ws.Range("H98").Value = 938.5714
$ws.Range("I98").Value = 924.4
$ws.Range("J98").Value = 974
$ws.Range("K98").Value = 924.4
$ws.Range("L98").Value = 974
$ws.Range("M98").Value = 573.6
$ws.Range("N98").Value = -3970
# Row 122
$ws.Range("H122").Value = 938.5714
$ws.Range("I122").Value = 924.4
$ws.Range("J122").Value = 974
$ws.Range("K122").Value = 2773.2
$ws.Range("L122").Value = 2922
$ws.Range("M122").Value = -323.1999999999998
$ws.Range("N122").Value = -7822
# Row 132
$ws.Range("H132").Value = 2815.7727
$ws.Range("I132").Value = 2668.1875
$ws.Range("K132").Value = 8004.5625
$ws.Range("M132").Value = -5474.5625
# Row 138
$ws.Range("H138").Value = 2595.7058
$ws.Range("I138").Value = 3147.2856
$ws.Range("J138").Value = 2209.6
$ws.Range("K138").Value = 9441.856800000001
$ws.Range("L138").Value = 6628.799999999999
$ws.Range("M138").Value = -4301.856800000001
$ws.Range("N138").Value = -16908.8

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3841.1482
$ws.Range("I32").Value = 2442.2727
$ws.Range("K32").Value = 2442.2727
$ws.Range("M32").Value = -2155.2727
# Row 45
$ws.Range("H45").Value = 2842.5334
$ws.Range("I45").Value = 1348.2858
$ws.Range("K45").Value = 1348.2858
$ws.Range("M45").Value = -971.2858000000001
# Row 61
$ws.Range("H61").Value = 1103.6666
$ws.Range("I61").Value = 905.5
$ws.Range("K61").Value = 905.5
$ws.Range("M61").Value = -693.5
# Row 63
$ws.Range("H63").Value = 4146.1
$ws.Range("I63").Value = 2280.1428
$ws.Range("K63").Value = 2280.1428
$ws.Range("M63").Value = -1594.1428
# Row 66
$ws.Range("H66").Value = 4146.1
$ws.Range("I66").Value = 2280.1428
$ws.Range("K66").Value = 11400.714
$ws.Range("M66").Value = -7968.714
# Row 122
$ws.Range("H122").Value = 2911.389
$ws.Range("I122").Value = 2640.6155
$ws.Range("J122").Value = 3615.4
$ws.Range("K122").Value = 7921.8465
$ws.Range("L122").Value = 10846.2
$ws.Range("M122").Value = -5471.8465
$ws.Range("N122").Value = -15746.2
# Row 132
$ws.Range("H132").Value = 1251.1666
$ws.Range("I132").Value = 1251.1666
$ws.Range("K132").Value = 3753.4998
$ws.Range("M132").Value = -1223.4998
# Row 136
$ws.Range("H136").Value = 1103.6666
$ws.Range("I136").Value = 905.5
$ws.Range("K136").Value = 2716.5
$ws.Range("M136").Value = -166.5

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 4266
$ws.Range("I105").Value = 4040.9167
$ws.Range("K105").Value = 4040.9167
$ws.Range("M105").Value = -2293.9167

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2122.4
$ws.Range("I31").Value = 2049.5
$ws.Range("K31").Value = 2049.5
$ws.Range("M31").Value = -1754.5
# Row 34
$ws.Range("H34").Value = 2122.4
$ws.Range("I34").Value = 2049.5
$ws.Range("K34").Value = 2049.5
$ws.Range("M34").Value = -1847.5
# Row 95
$ws.Range("H95").Value = 18690.666
$ws.Range("J95").Value = 18690.666
$ws.Range("L95").Value = 18690.666
$ws.Range("N95").Value = -24182.666
# Row 107
$ws.Range("H107").Value = 700
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
# Row 122
$ws.Range("H122").Value = 896.5833
$ws.Range("I122").Value = 840.3333
$ws.Range("J122").Value = 952.8333
$ws.Range("K122").Value = 2520.9999
$ws.Range("L122").Value = 2858.4999
$ws.Range("M122").Value = -70.9998999999998
$ws.Range("N122").Value = -7758.4999
# Row 132
$ws.Range("H132").Value = 5723.25
$ws.Range("I132").Value = 6398
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 19194
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -16664
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("CUL")
# Row 64
$ws.Range("H64").Value = 3000
$ws.Range("J64").Value = 3000
$ws.Range("L64").Value = 9000
$ws.Range("N64").Value = -9540
# Row 67
$ws.Range("H67").Value = 3000
$ws.Range("J67").Value = 3000
$ws.Range("L67").Value = 9000
$ws.Range("N67").Value = -10872
# Row 92
$ws.Range("H92").Value = 484
$ws.Range("I92").Value = 288.5
$ws.Range("J92").Value = 875
$ws.Range("K92").Value = 865.5
$ws.Range("L92").Value = 2625
$ws.Range("M92").Value = 382.5
$ws.Range("N92").Value = -5121
# Row 132
$ws.Range("H132").Value = 1699.5
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 832.25
$ws.Range("I97").Value = 531.25
$ws.Range("K97").Value = 531.25
$ws.Range("M97").Value = -35.25
# Row 104
$ws.Range("H104").Value = 5635.5
$ws.Range("J104").Value = 5635.5
$ws.Range("L104").Value = 5635.5
$ws.Range("N104").Value = -12623.5
# Row 105
$ws.Range("H105").Value = 9000
$ws.Range("J105").Value = 9000
$ws.Range("L105").Value = 9000
$ws.Range("N105").Value = -15988
# Row 132
$ws.Range("H132").Value = 10000
$ws.Range("I132").Value = 10000
$ws.Range("K132").Value = 30000
$ws.Range("M132").Value = -27470

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1807.4615
$ws.Range("I46").Value = 1461.5385
$ws.Range("J46").Value = 2153.3845
$ws.Range("K46").Value = 1461.5385
$ws.Range("L46").Value = 2153.3845
$ws.Range("M46").Value = -1273.5385
$ws.Range("N46").Value = -2529.3845
# Row 55
$ws.Range("H55").Value = 282.08334
$ws.Range("J55").Value = 298.1
$ws.Range("L55").Value = 298.1
$ws.Range("N55").Value = -644.1
# Row 61
$ws.Range("H61").Value = 1200
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 1200
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 1200
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -1604
# Row 108
$ws.Range("H108").Value = 40625.5
$ws.Range("J108").Value = 40625.5
$ws.Range("L108").Value = 40625.5
$ws.Range("N108").Value = -48305.5
# Row 110
$ws.Range("H110").Value = 50644
$ws.Range("J110").Value = 50644
$ws.Range("L110").Value = 50644
$ws.Range("N110").Value = -58824
# Row 112
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
# Row 113
$ws.Range("H113").Value = 1200
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1200
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5540

$ws = $wb.Worksheets.Item("WVR")
# Row 98
$ws.Range("H98").Value = 10590
$ws.Range("J98").Value = 10590
$ws.Range("L98").Value = 10590
$ws.Range("N98").Value = -16580
# Row 100
$ws.Range("H100").Value = 6338014.5
$ws.Range("I100").Value = 9957767
$ws.Range("K100").Value = 19915534
$ws.Range("M100").Value = -19914993
# Row 105
$ws.Range("H105").Value = 41091.777
$ws.Range("J105").Value = 41091.777
$ws.Range("L105").Value = 41091.777
$ws.Range("N105").Value = -48079.777
# Row 107
$ws.Range("H107").Value = 478.8
$ws.Range("J107").Value = 799
$ws.Range("L107").Value = 2397
$ws.Range("N107").Value = -6237
# Row 113
$ws.Range("H113").Value = 752.3333
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 4500
$ws.Range("N113").Value = -8840
# Row 126
$ws.Range("H126").Value = 3798.8
$ws.Range("I126").Value = 3798.8
$ws.Range("K126").Value = 11396.4
$ws.Range("M126").Value = -8926.400000000001

Write-Host "Edit complete"